$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.557.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.840.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4248'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3665'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07239'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8697'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.79'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.847.82'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.397'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.522'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06932'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.12'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009004'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.49'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.654.52'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.056'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.86'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.118.79'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.955'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.03'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.37'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.259'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.31'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.842'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08877'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7746'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.564'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.954'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.152'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.099'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05381'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.834'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5124'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1661'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.81%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.772'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.503'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.48'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06539'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.06'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4713'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.002'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.632'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.797'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.38%  '
